$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")
$ws.Rows("1:1").Select() | Out-Null
$ws.Rows("1:1").Delete()
$ws.Rows("1:1").Select() | Out-Null
$ws.Rows("1:1").Delete()
